$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = "Riccardo Briosi"
$ws.Range("B53").Value = "Riccardo Versini | Modium"
$ws.Range("C53").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("D53").Value = "Luca Frasca | Clitoriders"
$ws.Range("E53").Value = "Gianni Sala | FC SALAGIARDINI"
$ws.Range("F53").Value = "Christian Torboli | 4SINS"
